# Apply the Sprint Burndown Chart Template edits described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "Create Account" backlog item gets progress entries ---
$ws.Range("D6").Value = "Create Account"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 3
$ws.Range("I6").ClearContents()

# --- Row 8: add initial estimate + progress ---
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 1

# --- Row 13: new backlog item "Create Account" with estimate/progress ---
$ws.Range("D13").Value = "Create Account"
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = -1

# --- Row 17: new backlog item "Adasd" ---
$ws.Range("D17").Value = "Adasd"

# --- Row 26 totals: extend summation ranges from row 12 to row 20 ---
$ws.Range("E26").Formula = "=SUM(E6:E20)"
$ws.Range("F26:O26").Formula = "=E26-SUM(F6:F20)"

# --- Sheet view: scroll/zoom/selection changes ---
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K22").Select()

# --- Data validation changes ---
# XlDVType: xlValidateWholeNumber=1, xlValidateDecimal=2
# XlFormatConditionOperator: xlBetween=1
$ws.Range("S8").Validation.Delete()
$ws.Range("S8").Validation.Add(1, 1, 1, 0, 1000000)

$ws.Range("S9:S10").Validation.Delete()
$ws.Range("S9:S10").Validation.Add(1, 1, 1, -100, 1000000)

$ws.Range("E6:E25").Validation.Delete()
$ws.Range("F6:K25").Validation.Delete()
$ws.Range("L6:O11").Validation.Delete()
$ws.Range("O12").Validation.Delete()
$ws.Range("L12").Validation.Delete()
$ws.Range("L13:O13").Validation.Delete()
$ws.Range("L14").Validation.Delete()
$ws.Range("O14").Validation.Delete()
$ws.Range("L15:O15").Validation.Delete()
$ws.Range("L16").Validation.Delete()
$ws.Range("O16").Validation.Delete()
$ws.Range("L17:O17").Validation.Delete()
$ws.Range("L18").Validation.Delete()
$ws.Range("O18").Validation.Delete()
$ws.Range("L19:O19").Validation.Delete()
$ws.Range("L20").Validation.Delete()
$ws.Range("O20").Validation.Delete()
$ws.Range("L21:O21").Validation.Delete()
$ws.Range("L22").Validation.Delete()
$ws.Range("O22").Validation.Delete()
$ws.Range("L23:O23").Validation.Delete()
$ws.Range("L24").Validation.Delete()
$ws.Range("O24").Validation.Delete()
$ws.Range("L25:O25").Validation.Delete()

$ws.Range("E6:O27").Validation.Add(2, 1, 1, -24, 24)
